# Weekly refresh: insert two new price rows (latest week) above the
# existing data block for "Acelga" at Terminal Hortofrutícola Agro
# Chillán, pushing the rest of the rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current first data row (509),
# shifting rows 509:592 down to 511:594.
$ws.Rows("509:510").Insert()

# --- New row 509 (Primera) ---
$r = 509
$ws.Cells.Item($r, 1).Value  = 7
$ws.Cells.Item($r, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($r, 3).Value  = "Ñuble"
$ws.Cells.Item($r, 4).Value  = 45209
$ws.Cells.Item($r, 5).Value  = 16
$ws.Cells.Item($r, 6).Value  = 100112009
$ws.Cells.Item($r, 7).Value  = "Acelga"
$ws.Cells.Item($r, 8).Value  = "Sin especificar"
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 200
$ws.Cells.Item($r, 11).Value = 700
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 700
$ws.Cells.Item($r, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item($r, 15).Value = "Región de Ñuble"
$ws.Cells.Item($r, 16).Value = 700
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = "Hortaliza"

# --- New row 510 (Segunda) ---
$r = 510
$ws.Cells.Item($r, 1).Value  = 7
$ws.Cells.Item($r, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($r, 3).Value  = "Ñuble"
$ws.Cells.Item($r, 4).Value  = 45209
$ws.Cells.Item($r, 5).Value  = 16
$ws.Cells.Item($r, 6).Value  = 100112009
$ws.Cells.Item($r, 7).Value  = "Acelga"
$ws.Cells.Item($r, 8).Value  = "Sin especificar"
$ws.Cells.Item($r, 9).Value  = "Segunda"
$ws.Cells.Item($r, 10).Value = 150
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item($r, 15).Value = "Región de Ñuble"
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = "Hortaliza"
